$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add I0 and IF headers in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell H1 so the new cells
# match the bold/centered/bordered header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Add the corresponding data values in row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
